$wb = $excel.ActiveWorkbook

# Add the new "logical functions" worksheet after the last existing sheet.
# (This also clears the "tabSelected" flag on the previously active sheet,
# ByteValue, since a workbook has only one selected tab at a time.)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "logical functions"

# Populate the two new cells with the "anyTrue" method text.
$newSheet.Range("C7").Value = "Method boolean checkOr()"
$newSheet.Range("C8").Value = "return anyTrue(new boolean[]{true, false});"

# Match the column width used in the authored sheet as closely as this host
# can express: the target OOXML column width is 34.5703125 characters, and
# this host's ColumnWidth setter quantizes to 1/6-character steps, so feed it
# the character width (~33.667) whose quantized result lands on 34.5 - the
# closest reachable value to 34.5703125.
$newSheet.Columns.Item(3).ColumnWidth = 33.6666667

# Make the new sheet the active tab / selection, matching the diff.
$newSheet.Range("C8").Select() | Out-Null
$newSheet.Activate() | Out-Null
